$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 1000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 1000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -532
$ws.Range("H23").Value = 1000
$ws.Range("I23").Value = 1000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 1000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -766
$ws.Range("H41").Value = 330.9
$ws.Range("I41").Value = 146.7
$ws.Range("J41").Value = 515.1
$ws.Range("K41").Value = 146.7
$ws.Range("L41").Value = 515.1
$ws.Range("M41").Value = 293.3
$ws.Range("N41").Value = -1395.1
$ws.Range("H86").Value = 5475.125
$ws.Range("I86").Value = 4520.4
$ws.Range("J86").Value = 5909.091
$ws.Range("K86").Value = 4520.4
$ws.Range("L86").Value = 5909.091
$ws.Range("M86").Value = -3397.4
$ws.Range("N86").Value = -8155.091
$ws.Range("H89").Value = 5475.125
$ws.Range("I89").Value = 4520.4
$ws.Range("J89").Value = 5909.091
$ws.Range("K89").Value = 22602
$ws.Range("L89").Value = 29545.455
$ws.Range("M89").Value = -16986
$ws.Range("N89").Value = -40777.455
$ws.Range("H126").Value = 83249.75
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 83249.75
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 83249.75
$ws.Range("N126").Value = -93129.75
$ws.Range("H132").Value = 937.13336
$ws.Range("I132").Value = 807.6429000000001
$ws.Range("J132").Value = 2750
$ws.Range("K132").Value = 2422.9287
$ws.Range("L132").Value = 8250
$ws.Range("M132").Value = 107.0712999999996
$ws.Range("N132").Value = -13310
$ws.Range("H137").Value = 13160611
$ws.Range("I137").Value = 50001750
$ws.Range("J137").Value = 3061.3035
$ws.Range("K137").Value = 150005250
$ws.Range("L137").Value = 9183.9105
$ws.Range("M137").Value = -150002700
$ws.Range("N137").Value = -14283.9105
$ws.Range("H138").Value = 2847.2737
$ws.Range("I138").Value = 1824.2
$ws.Range("J138").Value = 3120.0933
$ws.Range("K138").Value = 5472.6
$ws.Range("L138").Value = 9360.2799
$ws.Range("M138").Value = -332.6000000000004
$ws.Range("N138").Value = -19640.2799
$ws.Range("H141").Value = 2730.8462
$ws.Range("I141").Value = 2730.8462
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8192.5386
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -3012.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4278.343
$ws.Range("I32").Value = 3238.3691
$ws.Range("J32").Value = 17798
$ws.Range("K32").Value = 3238.3691
$ws.Range("L32").Value = 17798
$ws.Range("M32").Value = -2951.3691
$ws.Range("N32").Value = -18372
$ws.Range("H39").Value = 5000
$ws.Range("I39").Value = 5000
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 5000
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -4480
$ws.Range("H45").Value = 90910790
$ws.Range("I45").Value = 100001720
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 100001720
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -100001343
$ws.Range("N45").Value = -2254
$ws.Range("H61").Value = 9289.286
$ws.Range("I61").Value = 6352.2
$ws.Range("J61").Value = 16632
$ws.Range("K61").Value = 6352.2
$ws.Range("L61").Value = 16632
$ws.Range("M61").Value = -6140.2
$ws.Range("N61").Value = -17056
$ws.Range("H117").Value = 12750
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 12750
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 12750
$ws.Range("N117").Value = -21928
$ws.Range("H136").Value = 9289.286
$ws.Range("I136").Value = 6352.2
$ws.Range("J136").Value = 16632
$ws.Range("K136").Value = 19056.6
$ws.Range("L136").Value = 49896
$ws.Range("M136").Value = -16506.6
$ws.Range("N136").Value = -54996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 69994
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 69994
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 69994
$ws.Range("N21").Value = -70466

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30409.422
$ws.Range("I31").Value = 2014.2609
$ws.Range("J31").Value = 73948.664
$ws.Range("K31").Value = 2014.2609
$ws.Range("L31").Value = 73948.664
$ws.Range("M31").Value = -1719.2609
$ws.Range("N31").Value = -74538.664
$ws.Range("H34").Value = 30409.422
$ws.Range("I34").Value = 2014.2609
$ws.Range("J34").Value = 73948.664
$ws.Range("K34").Value = 2014.2609
$ws.Range("L34").Value = 73948.664
$ws.Range("M34").Value = -1812.2609
$ws.Range("N34").Value = -74352.664
$ws.Range("H35").Value = 2161.5
$ws.Range("I35").Value = 2161.5
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 2161.5
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1867.5
$ws.Range("N35").ClearContents()
$ws.Range("H62").Value = 9822.556
$ws.Range("I62").Value = 4079.4
$ws.Range("J62").Value = 17001.5
$ws.Range("K62").Value = 4079.4
$ws.Range("L62").Value = 17001.5
$ws.Range("M62").Value = -3455.4
$ws.Range("N62").Value = -18249.5
$ws.Range("H65").Value = 9822.556
$ws.Range("I65").Value = 4079.4
$ws.Range("J65").Value = 17001.5
$ws.Range("K65").Value = 20397
$ws.Range("L65").Value = 85007.5
$ws.Range("M65").Value = -17277
$ws.Range("N65").Value = -91247.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 18502
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 18502
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 55506
$ws.Range("N42").Value = -56574
$ws.Range("H122").Value = 6019282
$ws.Range("I122").Value = 3663635.5
$ws.Range("J122").Value = 6804497.5
$ws.Range("K122").Value = 32972719.5
$ws.Range("L122").Value = 61240477.5
$ws.Range("M122").Value = -32970269.5
$ws.Range("N122").Value = -61245377.5
$ws.Range("H137").Value = 73891.36
$ws.Range("I137").Value = 1150
$ws.Range("J137").Value = 102987.9
$ws.Range("K137").Value = 3450
$ws.Range("L137").Value = 308963.7
$ws.Range("M137").Value = 1650
$ws.Range("N137").Value = -319163.7
$ws.Range("H138").Value = 2345.8572
$ws.Range("I138").Value = 1070.1666
$ws.Range("J138").Value = 10000
$ws.Range("K138").Value = 3210.4998
$ws.Range("L138").Value = 30000
$ws.Range("M138").Value = 1929.5002
$ws.Range("N138").Value = -40280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7911.625
$ws.Range("I7").Value = 5309.6665
$ws.Range("J7").Value = 11257
$ws.Range("K7").Value = 5309.6665
$ws.Range("L7").Value = 11257
$ws.Range("M7").Value = -5197.6665
$ws.Range("N7").Value = -11481
$ws.Range("H122").Value = 177530.61
$ws.Range("I122").Value = 238176.6
$ws.Range("J122").Value = 5700.3335
$ws.Range("K122").Value = 714529.8
$ws.Range("L122").Value = 17101.0005
$ws.Range("M122").Value = -712079.8
$ws.Range("N122").Value = -22001.0005
$ws.Range("H126").Value = 7911.625
$ws.Range("I126").Value = 5309.6665
$ws.Range("J126").Value = 11257
$ws.Range("K126").Value = 15928.9995
$ws.Range("L126").Value = 33771
$ws.Range("M126").Value = -13458.9995
$ws.Range("N126").Value = -38711
$ws.Range("H136").Value = 2713.5425
$ws.Range("I136").Value = 1719.8667
$ws.Range("J136").Value = 5907.5
$ws.Range("K136").Value = 5159.6001
$ws.Range("L136").Value = 17722.5
$ws.Range("M136").Value = -2609.6001
$ws.Range("N136").Value = -22822.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 20000
$ws.Range("I41").Value = 20000
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 20000
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -19610
$ws.Range("N41").ClearContents()
$ws.Range("H126").Value = 3114.2856
$ws.Range("I126").Value = 4455.75
$ws.Range("J126").Value = 1325.6666
$ws.Range("K126").Value = 13367.25
$ws.Range("L126").Value = 3976.9998
$ws.Range("M126").Value = -10897.25
$ws.Range("N126").Value = -8916.9998
$ws.Range("H132").Value = 2827.8948
$ws.Range("I132").Value = 1866.1765
$ws.Range("J132").Value = 11002.5
$ws.Range("K132").Value = 5598.529500000001
$ws.Range("L132").Value = 33007.5
$ws.Range("M132").Value = -3068.529500000001
$ws.Range("N132").Value = -38067.5
$ws.Range("H136").Value = 3100.4614
$ws.Range("I136").Value = 1966.0256
$ws.Range("J136").Value = 6503.769
$ws.Range("K136").Value = 5898.0768
$ws.Range("L136").Value = 19511.307
$ws.Range("M136").Value = -3348.0768
$ws.Range("N136").Value = -24611.307
